$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "29.174.79"
$cell.Style = "Normal"

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "1.825.51"
$cell.Style = "Normal"

$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = "0.9987"
$cell.Style = "Normal"
$ws.Range("E4").Value = "  -0.01%  "

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "234.86"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  -1.81%  "

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "0.6002"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  -4.07%  "

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.07062"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  -4.86%  "

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.2792"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  -3.53%  "

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "23.50"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  -5.12%  "

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.07604"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  -1.49%  "

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "1.829.90"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  -0.51%  "

$ws.Range("E13").Value = "  -3.21%  "

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "0.6300"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  -6.44%  "

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "0.000009908"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  -2.84%  "

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "2.076.89"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  -0.51%  "

$ws.Range("E17").Value = "  -3.50%  "

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "5.858"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  -5.81%  "

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "29.170.62"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  -0.72%  "

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "226.58"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  -2.63%  "

$ws.Range("E21").Value = "  +0.03%  "

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "11.72"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  -4.61%  "

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "6.990"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  -4.73%  "

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "1.001"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  +0.10%  "

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "154.85"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  -2.05%  "

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "8.016"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  -5.31%  "

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "0.1298"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  -3.45%  "

$ws.Range("E28").Value = "  -4.46%  "

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "1.491"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  +2.54%  "

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "0.06235"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  -14.24%  "

$ws.Range("E31").Value = "  -1.74%  "

$ws.Range("E32").Value = "  -5.08%  "

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "3.797"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  -5.76%  "

$ws.Range("E34").Value = "  -1.44%  "

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "1.738"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  -4.30%  "

$ws.Range("E36").Value = "  -7.83%  "

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "2.532"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  -1.51%  "

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "1.216.29"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  -1.08%  "

$ws.Range("E39").Value = "  -3.03%  "

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "0.01734"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  -5.31%  "

$ws.Range("E41").Value = "  -5.87%  "

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "0.9055"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  -3.89%  "

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "0.9995"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  +0.00%  "

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "1.982.83"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  -0.33%  "

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "100.31"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  -0.19%  "

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "62.74"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  -4.10%  "

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "0.00000000117"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  -2.34%  "

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "8.528"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  -3.79%  "

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "1.597"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  -6.17%  "

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "0.4548"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  -0.74%  "

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "0.05500"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  -2.69%  "
